# Add a new "2021" data column (N) to the 6.1.1 water-services indicator sheet.
# Mirrors the author's manual edit: a new year column was appended after the
# existing M (2020) column, rows 4-14 got their existing/new "164" numeric
# format, row 4 (Kyrgyz Republic total) became bold, and the selection cursor
# was left on the newly entered N2 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 2 (thin separator row, border only, no value) -------------------
$ws.Range("D2").Copy()
$ws.Range("N2").PasteSpecial($xlPasteFormats)

# --- Row 3 (year headers) --------------------------------------------------
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial($xlPasteFormats)
$ws.Range("N3").Value = 2021

# --- Row 4 (Kyrgyz Republic total, bold "164"-format row) -----------------
# D4:L4 pick up the same bold style M4 already used.
$ws.Range("D4:L4").Font.Bold = $true

$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial($xlPasteFormats)
$ws.Range("N4").Value = 95.134712433469176

# --- Rows 5-14 (oblast / city data rows, plain "164" numeric format) ------
$ws.Range("D5").Copy()
$ws.Range("N5").PasteSpecial($xlPasteFormats)
$ws.Range("N5").Value = 99.705541665880986

$ws.Range("D6").Copy()
$ws.Range("N6").PasteSpecial($xlPasteFormats)
$ws.Range("N6").Value = 92.425193326577897

$ws.Range("D7").Copy()
$ws.Range("N7").PasteSpecial($xlPasteFormats)
$ws.Range("N7").Value = 88.209991167538519

$ws.Range("D8").Copy()
$ws.Range("N8").PasteSpecial($xlPasteFormats)
$ws.Range("N8").Value = 92.225038985690773

$ws.Range("D9").Copy()
$ws.Range("N9").PasteSpecial($xlPasteFormats)
$ws.Range("N9").Value = 96.801032063987265

$ws.Range("D10").Copy()
$ws.Range("N10").PasteSpecial($xlPasteFormats)
$ws.Range("N10").Value = 97.660491031729507

$ws.Range("D11").Copy()
$ws.Range("N11").PasteSpecial($xlPasteFormats)
$ws.Range("N11").Value = 90.23262877800066

$ws.Range("D12").Copy()
$ws.Range("N12").PasteSpecial($xlPasteFormats)
$ws.Range("N12").Value = 99.653994395099105

$ws.Range("D13").Copy()
$ws.Range("N13").PasteSpecial($xlPasteFormats)
$ws.Range("N13").Value = 100

$ws.Range("D14").Copy()
$ws.Range("N14").PasteSpecial($xlPasteFormats)
$ws.Range("N14").Value = 100

# --- Row 15 (bottom border, right-aligned "164" numeric format) -----------
$ws.Range("M15").Copy()
$ws.Range("N15").PasteSpecial($xlPasteFormats)
$ws.Range("N15").Value = 100

# --- Misc sheet-level bits --------------------------------------------------
# Bump the print resolution the same way the source workbook recorded it.
$ws.PageSetup.PrintQuality = 300

# Leave the selection where the author's last edit landed: the freshly typed N2.
$ws.Range("N2").Select()
